{"js": "// Replace each \"AxB=\" multiplication problem in the document's table\n// with its new value, per the commit diff. All old values are unique\n// in the document, so a straightforward search-and-replace per pair\n// is unambiguous and preserves run formatting (font/size) because we\n// replace the text inside the matched Range rather than rebuilding runs.\nconst replacements = [\n  [\"26\u00d729=\", \"96\u00d757=\"],\n  [\"45\u00d711=\", \"13\u00d786=\"],\n  [\"40\u00d782=\", \"56\u00d778=\"],\n  [\"88\u00d777=\", \"56\u00d723=\"],\n  [\"54\u00d751=\", \"37\u00d773=\"],\n  [\"47\u00d729=\", \"33\u00d760=\"],\n  [\"79\u00d781=\", \"83\u00d740=\"],\n  [\"83\u00d783=\", \"84\u00d727=\"],\n  [\"60\u00d778=\", \"77\u00d737=\"],\n  [\"25\u00d771=\", \"89\u00d758=\"],\n  [\"78\u00d739=\", \"12\u00d788=\"],\n  [\"51\u00d798=\", \"39\u00d713=\"],\n  [\"46\u00d725=\", \"24\u00d724=\"],\n  [\"16\u00d768=\", \"46\u00d736=\"],\n  [\"87\u00d714=\", \"99\u00d738=\"],\n  [\"40\u00d740=\", \"25\u00d754=\"],\n  [\"55\u00d750=\", \"73\u00d746=\"],\n  [\"19\u00d770=\", \"26\u00d764=\"],\n  [\"36\u00d741=\", \"87\u00d741=\"],\n  [\"77\u00d746=\", \"14\u00d771=\"],\n  [\"37\u00d772=\", \"22\u00d729=\"],\n  [\"23\u00d711=\", \"74\u00d790=\"],\n  [\"95\u00d788=\", \"67\u00d794=\"],\n  [\"77\u00d778=\", \"27\u00d729=\"],\n  [\"78\u00d766=\", \"14\u00d765=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"AxB=\" multiplication problem in the document's table\n# with its new value, per the commit diff. All old values are unique\n# in the document, so Find/Replace (wdReplaceAll) per pair is\n# unambiguous and preserves run formatting (font/size) since Word's\n# Find.Execute replace-in-place keeps the surrounding run properties.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"26\u00d729=\"; New = \"96\u00d757=\" },\n    @{ Old = \"45\u00d711=\"; New = \"13\u00d786=\" },\n    @{ Old = \"40\u00d782=\"; New = \"56\u00d778=\" },\n    @{ Old = \"88\u00d777=\"; New = \"56\u00d723=\" },\n    @{ Old = \"54\u00d751=\"; New = \"37\u00d773=\" },\n    @{ Old = \"47\u00d729=\"; New = \"33\u00d760=\" },\n    @{ Old = \"79\u00d781=\"; New = \"83\u00d740=\" },\n    @{ Old = \"83\u00d783=\"; New = \"84\u00d727=\" },\n    @{ Old = \"60\u00d778=\"; New = \"77\u00d737=\" },\n    @{ Old = \"25\u00d771=\"; New = \"89\u00d758=\" },\n    @{ Old = \"78\u00d739=\"; New = \"12\u00d788=\" },\n    @{ Old = \"51\u00d798=\"; New = \"39\u00d713=\" },\n    @{ Old = \"46\u00d725=\"; New = \"24\u00d724=\" },\n    @{ Old = \"16\u00d768=\"; New = \"46\u00d736=\" },\n    @{ Old = \"87\u00d714=\"; New = \"99\u00d738=\" },\n    @{ Old = \"40\u00d740=\"; New = \"25\u00d754=\" },\n    @{ Old = \"55\u00d750=\"; New = \"73\u00d746=\" },\n    @{ Old = \"19\u00d770=\"; New = \"26\u00d764=\" },\n    @{ Old = \"36\u00d741=\"; New = \"87\u00d741=\" },\n    @{ Old = \"77\u00d746=\"; New = \"14\u00d771=\" },\n    @{ Old = \"37\u00d772=\"; New = \"22\u00d729=\" },\n    @{ Old = \"23\u00d711=\"; New = \"74\u00d790=\" },\n    @{ Old = \"95\u00d788=\"; New = \"67\u00d794=\" },\n    @{ Old = \"77\u00d778=\"; New = \"27\u00d729=\" },\n    @{ Old = \"78\u00d766=\"; New = \"14\u00d765=\" }\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Execute($pair.Old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.New, \"wdReplaceAll\")\n}\n"}
